$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4,D5,D6,D7,D9,D10,D11,D12,D14,D15,D16,D18,D19,D20,D22,D24,D25,D26,D27,D28,D29,D30,D31,D32,D33,D34,D35,D36,D37,D38,D39,D40,D41,D42,D43,D44,D45,D46,D48,D49,D50,D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.306.24'
$ws.Range("E2").Value = '  +0.29%  '
$ws.Range("D3").Value = '1.865.72'
$ws.Range("E3").Value = '  +0.28%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '234.31'
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("D7").Value = '0.4695'
$ws.Range("E7").Value = '  +0.32%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = '0.06568'
$ws.Range("E9").Value = '  +0.31%  '
$ws.Range("D10").Value = '21.63'
$ws.Range("E10").Value = '  -0.55%  '
$ws.Range("D11").Value = '0.07876'
$ws.Range("E11").Value = '  -0.64%  '
$ws.Range("D12").Value = '96.41'
$ws.Range("E12").Value = '  -1.32%  '
$ws.Range("D13").Value = '1.866.53'
$ws.Range("E13").Value = '  +0.24%  '
$ws.Range("D14").Value = '0.6935'
$ws.Range("E14").Value = '  +1.76%  '
$ws.Range("D15").Value = '5.103'
$ws.Range("E15").Value = '  -1.44%  '
$ws.Range("D16").Value = '267.47'
$ws.Range("E16").Value = '  -0.30%  '
$ws.Range("D17").Value = '30.242.60'
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("D18").Value = '14.00'
$ws.Range("E18").Value = '  +1.66%  '
$ws.Range("D19").Value = '0.000007674'
$ws.Range("E19").Value = '  +2.92%  '
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.11%  '
$ws.Range("D21").Value = '2.111.35'
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").Value = '1.002'
$ws.Range("E23").Value = '  -1.83%  '
$ws.Range("D24").Value = '6.190'
$ws.Range("D25").Value = '9.382'
$ws.Range("E25").Value = '  +1.70%  '
$ws.Range("D26").Value = '167.24'
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("D27").Value = '18.81'
$ws.Range("E27").Value = '  -0.59%  '
$ws.Range("D28").Value = '1.941'
$ws.Range("E28").Value = '  -1.03%  '
$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D29").Value = '0.09873'
$ws.Range("E29").Value = '  +0.18%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '1.353'
$ws.Range("E30").Value = '  -2.36%  '
$ws.Range("D31").Value = '4.368'
$ws.Range("E31").Value = '  -0.51%  '
$ws.Range("D32").Value = '1.458'
$ws.Range("E32").Value = '  -1.12%  '
$ws.Range("D33").Value = '4.065'
$ws.Range("E33").Value = '  -0.19%  '
$ws.Range("D34").Value = '0.04756'
$ws.Range("E34").Value = '  +0.87%  '
$ws.Range("D35").Value = '1.134'
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("D36").Value = '0.7015'
$ws.Range("E36").Value = '  -0.32%  '
$ws.Range("D37").Value = '2.725'
$ws.Range("E37").Value = '  +0.70%  '
$ws.Range("D38").Value = '0.01872'
$ws.Range("E38").Value = '  -0.34%  '
$ws.Range("D39").Value = '2.799'
$ws.Range("E39").Value = '  +6.95%  '
$ws.Range("D40").Value = '6.230'
$ws.Range("E40").Value = '  -0.32%  '
$ws.Range("D41").Value = '72.91'
$ws.Range("E41").Value = '  -2.15%  '
$ws.Range("D42").Value = '1.947'
$ws.Range("E42").Value = '  +0.18%  '
$ws.Range("D43").Value = '0.8422'
$ws.Range("E43").Value = '  -0.45%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '0.4171'
$ws.Range("E44").Value = '  -0.08%  '
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").Value = '  +0.17%  '
$ws.Range("D46").Value = '102.33'
$ws.Range("E46").Value = '  -0.91%  '
$ws.Range("D48").Value = '942.30'
$ws.Range("E48").Value = '  -1.43%  '
$ws.Range("D49").Value = '9.066'
$ws.Range("E49").Value = '  -1.77%  '
$ws.Range("D50").Value = '34.48'
$ws.Range("E50").Value = '  +0.93%  '
$ws.Range("D51").Value = '0.05674'
$ws.Range("E51").Value = '  +0.22%  '
